{"js": "// Update the division problems in the practice-sheet table.\n// Each cell holds a single \"NNN\u00f7N=\" run; replace the dividend/divisor\n// pairs with the regenerated values while keeping formatting untouched.\n\nconst replacements = [\n  { find: \"382\u00f72=\", replace: \"485\u00f77=\" },\n  { find: \"239\u00f79=\", replace: \"519\u00f75=\" },\n  { find: \"447\u00f77=\", replace: \"666\u00f76=\" },\n  { find: \"787\u00f74=\", replace: \"952\u00f76=\" },\n  { find: \"960\u00f79=\", replace: \"224\u00f75=\" },\n  { find: \"552\u00f76=\", replace: \"964\u00f79=\" },\n  { find: \"662\u00f79=\", replace: \"919\u00f76=\" },\n  { find: \"161\u00f77=\", replace: \"683\u00f76=\" },\n  { find: \"443\u00f79=\", replace: \"309\u00f79=\" },\n  { find: \"361\u00f75=\", replace: \"971\u00f77=\" },\n  { find: \"372\u00f75=\", replace: \"900\u00f73=\" },\n  { find: \"643\u00f75=\", replace: \"555\u00f78=\" },\n  { find: \"452\u00f75=\", replace: \"236\u00f78=\" },\n  { find: \"877\u00f78=\", replace: \"972\u00f74=\" },\n  { find: \"204\u00f74=\", replace: \"731\u00f76=\" },\n  { find: \"910\u00f78=\", replace: \"277\u00f79=\" },\n  { find: \"910\u00f72=\", replace: \"156\u00f78=\" },\n  { find: \"718\u00f76=\", replace: \"511\u00f77=\" },\n  { find: \"404\u00f76=\", replace: \"332\u00f74=\" },\n  { find: \"871\u00f78=\", replace: \"884\u00f74=\" },\n  { find: \"835\u00f74=\", replace: \"450\u00f75=\" },\n  { find: \"193\u00f73=\", replace: \"205\u00f79=\" },\n  { find: \"491\u00f72=\", replace: \"176\u00f77=\" },\n  { find: \"348\u00f75=\", replace: \"874\u00f76=\" },\n  { find: \"915\u00f79=\", replace: \"129\u00f75=\" },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the division problems in the practice-sheet table.\n# Each cell holds a single \"NNN\u00f7N=\" run; replace the dividend/divisor\n# pairs with the regenerated values while keeping formatting untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"382\u00f72=\"; New = \"485\u00f77=\" },\n    @{ Old = \"239\u00f79=\"; New = \"519\u00f75=\" },\n    @{ Old = \"447\u00f77=\"; New = \"666\u00f76=\" },\n    @{ Old = \"787\u00f74=\"; New = \"952\u00f76=\" },\n    @{ Old = \"960\u00f79=\"; New = \"224\u00f75=\" },\n    @{ Old = \"552\u00f76=\"; New = \"964\u00f79=\" },\n    @{ Old = \"662\u00f79=\"; New = \"919\u00f76=\" },\n    @{ Old = \"161\u00f77=\"; New = \"683\u00f76=\" },\n    @{ Old = \"443\u00f79=\"; New = \"309\u00f79=\" },\n    @{ Old = \"361\u00f75=\"; New = \"971\u00f77=\" },\n    @{ Old = \"372\u00f75=\"; New = \"900\u00f73=\" },\n    @{ Old = \"643\u00f75=\"; New = \"555\u00f78=\" },\n    @{ Old = \"452\u00f75=\"; New = \"236\u00f78=\" },\n    @{ Old = \"877\u00f78=\"; New = \"972\u00f74=\" },\n    @{ Old = \"204\u00f74=\"; New = \"731\u00f76=\" },\n    @{ Old = \"910\u00f78=\"; New = \"277\u00f79=\" },\n    @{ Old = \"910\u00f72=\"; New = \"156\u00f78=\" },\n    @{ Old = \"718\u00f76=\"; New = \"511\u00f77=\" },\n    @{ Old = \"404\u00f76=\"; New = \"332\u00f74=\" },\n    @{ Old = \"871\u00f78=\"; New = \"884\u00f74=\" },\n    @{ Old = \"835\u00f74=\"; New = \"450\u00f75=\" },\n    @{ Old = \"193\u00f73=\"; New = \"205\u00f79=\" },\n    @{ Old = \"491\u00f72=\"; New = \"176\u00f77=\" },\n    @{ Old = \"348\u00f75=\"; New = \"874\u00f76=\" },\n    @{ Old = \"915\u00f79=\"; New = \"129\u00f75=\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2) | Out-Null\n}\n\n$d.Save()\n"}
